$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data (data rows shift from 1-11 to 2-12)
[void]$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Coluna 1"
$ws.Range("B1").Value = "Coluna 2"
$ws.Range("C1").Value = "Coluna 3"
$ws.Range("D1").Value = "Coluna 4"
$ws.Range("E1").Value = "Coluna 5"
$ws.Range("F1").Value = "Coluna 6"
$ws.Range("G1").Value = "Coluna 7"
$ws.Range("H1").Value = "Coluna 8"
$ws.Range("I1").Value = "Coluna 9"
$ws.Range("J1").Value = "Coluna 10"
$ws.Range("K1").Value = "Coluna 11"

# Narrow a few columns that used to be wider than needed
$ws.Columns.Item(2).ColumnWidth = 23.5
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws.Columns.Item(7).ColumnWidth = 21.666666666666668
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668

# Select the new header row, matching the post-edit selection
[void]$ws.Range("A1:K1").Select()
